# Test-Generator-GenerateUploadSampleData.xlsx
#
# Commit: "New rule: All functional test data of all types must be found
# by q=__TEST__. This will make it easier to clean up later."
#
# 1) The Transaction sheet's sample rows are regenerated with new dates
#    (and the "Big Dogs" $4000 transaction moves from row 14 to row 4).
# 2) Every occurrence of the sample Payee names / Categories is renamed
#    so they can be found with the search term "__TEST__":
#       Big Money -> AA__TEST__1      X:Y -> AA:__TEST__:A
#       Big Stuff -> AA__TEST__2      X:Z -> AA:__TEST__:B
#       Big Dogs  -> AA__TEST__3      A:B -> AA:__TEST__:C
#                                     C:D -> AA:__TEST__:D
#    These names/categories are reused verbatim on the Payee, Split and
#    BudgetTx sheets, so renaming them there too keeps everything in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transaction")

# --- Row-by-row data fixups (still using the OLD labels; the rename pass
#     at the end of this script takes care of the text itself) ---

$ws.Range("B2").Value = 44214
$ws.Range("B3").Value = 44223

$ws.Range("B4").Value = 44230
$ws.Range("C4").Value = "Big Dogs"
$ws.Range("D4").Value = -4000
$ws.Range("E4").Value = "C:D"

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 44245
$ws.Range("C5").Value = "Big Stuff"
$ws.Range("E5").Value = "A:B"
$ws.Range("M5").Value = $false

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 44254
$ws.Range("C6").Value = "Big Money"
$ws.Range("E6").ClearContents()
$ws.Range("M6").Value = $true

$ws.Range("A7").Value = 0
$ws.Range("B7").Value = 44273
$ws.Range("C7").Value = "Big Stuff"
$ws.Range("E7").Value = "A:B"
$ws.Range("M7").Value = $false

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 44282
$ws.Range("C8").Value = "Big Money"
$ws.Range("E8").ClearContents()
$ws.Range("M8").Value = $true

$ws.Range("A9").Value = 0
$ws.Range("B9").Value = 44304
$ws.Range("C9").Value = "Big Stuff"
$ws.Range("E9").Value = "A:B"
$ws.Range("M9").Value = $false

$ws.Range("A10").Value = 4
$ws.Range("B10").Value = 44313
$ws.Range("C10").Value = "Big Money"
$ws.Range("E10").ClearContents()
$ws.Range("M10").Value = $true

$ws.Range("A11").Value = 0
$ws.Range("B11").Value = 44334
$ws.Range("C11").Value = "Big Stuff"
$ws.Range("E11").Value = "A:B"
$ws.Range("M11").Value = $false

$ws.Range("A12").Value = 5
$ws.Range("B12").Value = 44343
$ws.Range("C12").Value = "Big Money"
$ws.Range("E12").ClearContents()
$ws.Range("M12").Value = $true

$ws.Range("A13").Value = 0
$ws.Range("B13").Value = 44365
$ws.Range("C13").Value = "Big Stuff"
$ws.Range("E13").Value = "A:B"
$ws.Range("M13").Value = $false

$ws.Range("A14").Value = 6
$ws.Range("B14").Value = 44374
$ws.Range("C14").Value = "Big Money"
$ws.Range("D14").Value = -250
$ws.Range("E14").ClearContents()
$ws.Range("M14").Value = $true

$ws.Range("B15").Value = 44395
$ws.Range("B16").Value = 44404
$ws.Range("B17").Value = 44426
$ws.Range("B18").Value = 44435
$ws.Range("B19").Value = 44457
$ws.Range("B20").Value = 44466
$ws.Range("B21").Value = 44487
$ws.Range("B22").Value = 44496
$ws.Range("B23").Value = 44518
$ws.Range("B24").Value = 44527
$ws.Range("B25").Value = 44548
$ws.Range("B26").Value = 44557

# --- Global rename pass: apply the new "__TEST__" naming convention to
#     every sheet that references these Payee names / Categories ---

$renames = @{
    "Big Stuff" = "AA__TEST__2"
    "Big Money" = "AA__TEST__1"
    "Big Dogs"  = "AA__TEST__3"
    "A:B"       = "AA:__TEST__:C"
    "C:D"       = "AA:__TEST__:D"
    "X:Y"       = "AA:__TEST__:A"
    "X:Z"       = "AA:__TEST__:B"
}

foreach ($sheet in $wb.Worksheets) {
    foreach ($old in $renames.Keys) {
        $new = $renames[$old]
        $sheet.Cells.Replace($old, $new)
    }
}
